# Fix import other request
#
# The "import" sheet had a spurious "Ghi chú" (Note) column at J that is not
# wired up to anything downstream; this removes it entirely so the form
# lines up with the expected import layout again. Removing the whole column
# shifts every column from K onward one slot to the left (K->J, L->K, ...,
# S->R), which Excel's "Delete Entire Column" takes care of for the cell
# grid, column widths, the sheet dimension and the data-validation range.
#
# Cell comments are anchored objects that this host does not re-home
# automatically when a column is deleted, so we snapshot their text first,
# delete them, perform the column delete, and then re-create them one
# column to the left of where they used to be.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("import")

# --- 1. Snapshot the comments that live to the right of the doomed column J
#        (comments on C1:I1 stay put and are left untouched). ---
$commentRefs = @("K1", "N1", "O1", "P1", "Q1", "R1", "S1")
$commentTexts = @{}
foreach ($ref in $commentRefs) {
    $cmt = $ws.Range($ref).Comment
    if ($cmt -ne $null) {
        $commentTexts[$ref] = $cmt.Text()
        $cmt.Delete()
    }
}

# --- 2. Delete the whole "Ghi chú" column (J). Everything to its right
#        (values, styles, column widths, dimension, data validation sqref,
#        shared strings) shifts left by one column automatically. ---
$ws.Columns("J").Delete()

# --- 3. Re-create the comments one column to the left of their old spot. ---
$commentMap = @{
    "K1" = "J1"
    "N1" = "M1"
    "O1" = "N1"
    "P1" = "O1"
    "Q1" = "P1"
    "R1" = "Q1"
    "S1" = "R1"
}
foreach ($oldRef in $commentRefs) {
    if ($commentTexts.ContainsKey($oldRef)) {
        $newRef = $commentMap[$oldRef]
        $ws.Range($newRef).AddComment($commentTexts[$oldRef]) | Out-Null
    }
}

# --- 4. The active selection moved on to the next entry row (J11). ---
$ws.Range("J11").Select() | Out-Null
